# Sprint 2 Backlog - Burndown.xlsx
#
# Commit: "Fixed bug with viewing recipes" - the backlog item in row 7
# ("Fixed bug that caused Recipe details to not show for the first
# recipe...") is updated to record 2 actual hours worked (column D,
# "Actual Hours"). The burndown total in D22 (=SUM(D3:D21)) recalculates
# automatically from 7.25 -> 9.25 as a result.
#
# The workbook was also left with cell E7 as the active selection when it
# was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 ("Fixed bug with viewing recipes"): record 2 actual hours in
# column D. This was previously blank.
$ws.Range("D7").Value = 2

# Scroll the view roughly into position (column B at the left edge, row 1
# at the top) to mirror the saved window state, then leave E7 selected as
# the active cell, matching the saved sheetView/selection.
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("E7").Select()
